$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the stale cell selection stored in the sheet view (was J3) back to the default
[void]$ws.Range("A1").Select()

# BP-813: Affiliate Mapping for True Independent Stations
# Rename header columns: "Affiliation Mismatch Note" -> "IsTrueIND", "SalesGroupName" -> "RepFirm"
$ws.Range("H1").Value = "IsTrueIND"
$ws.Range("J1").Value = "RepFirm"

# Align formatting of the renamed header cells with the rest of the header row
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Align formatting of columns I and J (rows 2-5) with the rest of the data rows
$ws.Range("H2").Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)
$ws.Range("I3:J3").PasteSpecial(-4122)
$ws.Range("I4:J4").PasteSpecial(-4122)
$ws.Range("I5:J5").PasteSpecial(-4122)

$excel.CutCopyMode = 0
